# fix some issue when update product
# Row 2 of the "Data" sheet held one sample receipt line (Product ID,
# Product Name, Quantity, Lot ID, Location Code). The product identifier
# columns are no longer populated on this row, and the Lot ID / Location
# Code need refreshing to the latest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Product ID / Product Name: clear the values but keep the cells (and
# their default/"Normal" style) present on the row.
$ws.Range("A2").Value = ""
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = ""
$ws.Range("B2").Style = "Normal"

# Quantity (C2) stays as-is ("1").

# Lot ID and Location Code get updated to the latest values. Force text
# formatting so the purely-numeric-looking values (e.g. "413") are
# stored as text like the rest of this column, then drop back to the
# sheet's normal/default style (no custom number format lingers on the
# cell).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "413"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1712140573166"
$ws.Range("E2").Style = "Normal"
